# Auto-generated edit script for cryptos.xlsx update
# Updates Price (D) and Volume(1h) (E) columns, plus a row-38/39 swap
# (Aptos <-> VeChain) to match the refreshed data feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "21.503.32"
$ws.Range("E2").Value = "  -2.68%  "
$ws.Range("D3").Value = "1.529.98"
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("E5").Value = "  +0.16%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "288.37"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.36%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.3873"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -2.88%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3168"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -2.03%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "42.65"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -3.15%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.07148"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -2.47%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "1.068"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -1.55%  "
$ws.Range("E12").Value = "  +0.22%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "5.721"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +0.31%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "18.14"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -4.42%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "6.535"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.89%  "
$ws.Range("D16").Value = "1.535.69"
$ws.Range("E16").Value = "  -1.26%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.00001088"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -4.62%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.06607"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.12%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "83.35"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.65%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.094"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -3.51%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "15.38"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -2.32%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "10.79"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -4.52%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.372"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.44%  "
$ws.Range("D25").Value = "21.510.10"
$ws.Range("E25").Value = "  -2.70%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.372"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -3.10%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "148.87"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.04%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "18.32"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.83%  "
$ws.Range("E29").Value = "  -0.70%  "
$ws.Range("D30").Value = "1.708.69"
$ws.Range("E30").Value = "  -1.22%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "116.46"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -2.40%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "6.050"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +4.92%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.9525"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -6.18%  "
$ws.Range("E34").Value = "  -4.31%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "8.490"
$c.Style = "Normal"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "5.155"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.07%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.488"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -8.18%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.02196"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -3.59%  "
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "11.27"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +4.47%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.05878"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -4.26%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.2019"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -2.12%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.180"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -3.27%  "
$ws.Range("E43").Value = "  +0.19%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.5744"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -2.23%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "13.08"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.11%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "3.718"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -1.26%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.5550"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -1.24%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.892"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -1.12%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.157"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +1.28%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "115.42"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -2.96%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.06668"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -2.72%  "
